$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.487834453582764
$ws.Range("B1").Value = 2.765366315841675
$ws.Range("C1").Value = 6.87101411819458
$ws.Range("D1").Value = 1.747084140777588
$ws.Range("E1").Value = 0.8957867622375488
